$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.99924139758967
$ws.Cells.Item(2, 3).Value = 6.36166412509978
$ws.Cells.Item(2, 4).Value = 8.045919013763262
$ws.Cells.Item(2, 5).Value = 12.61410445800567
$ws.Cells.Item(2, 6).Value = 36.23237742641413
$ws.Cells.Item(2, 9).Value = 27.12140291937802
$ws.Cells.Item(2, 10).Value = 9.893607410861655
$ws.Cells.Item(2, 11).Value = 11.2597383138669
$ws.Cells.Item(2, 12).Value = 11.15898360159018
$ws.Cells.Item(2, 14).Value = 20.53014082243471
$ws.Cells.Item(2, 15).Value = 27.93921226684197
$ws.Cells.Item(3, 2).Value = 13.80182180899685
$ws.Cells.Item(3, 3).Value = 6.29821215188395
$ws.Cells.Item(3, 4).Value = 8.029415387163716
$ws.Cells.Item(3, 5).Value = 12.63012674708924
$ws.Cells.Item(3, 6).Value = 36.31096809774383
$ws.Cells.Item(3, 9).Value = 27.20787815117196
$ws.Cells.Item(3, 10).Value = 9.911393326257361
$ws.Cells.Item(3, 11).Value = 11.12130971306549
$ws.Cells.Item(3, 12).Value = 11.15848792935484
$ws.Cells.Item(3, 14).Value = 20.58736145625731
$ws.Cells.Item(3, 15).Value = 28.01810763834343
$ws.Cells.Item(4, 2).Value = 13.68148359534495
$ws.Cells.Item(4, 3).Value = 6.258371957166721
$ws.Cells.Item(4, 4).Value = 8.020352362113471
$ws.Cells.Item(4, 5).Value = 12.64144690360771
$ws.Cells.Item(4, 6).Value = 36.36595609200809
$ws.Cells.Item(4, 9).Value = 27.26507884983419
$ws.Cells.Item(4, 10).Value = 9.922986257986144
$ws.Cells.Item(4, 11).Value = 11.03714018041005
$ws.Cells.Item(4, 12).Value = 11.15955070620238
$ws.Cells.Item(4, 14).Value = 20.62416361594976
$ws.Cells.Item(4, 15).Value = 28.07121664278751
$ws.Cells.Item(5, 2).Value = 13.63272346479588
$ws.Cells.Item(5, 3).Value = 6.24192256064686
$ws.Cells.Item(5, 4).Value = 8.016930988503375
$ws.Cells.Item(5, 5).Value = 12.64643330601867
$ws.Cells.Item(5, 6).Value = 36.39005591659568
$ws.Cells.Item(5, 9).Value = 27.28942094007846
$ws.Cells.Item(5, 10).Value = 9.927879950347625
$ws.Cells.Item(5, 11).Value = 11.00308519917024
$ws.Cells.Item(5, 12).Value = 11.16032879750206
$ws.Cells.Item(5, 14).Value = 20.63958154309454
$ws.Cells.Item(5, 15).Value = 28.09403195708304
$ws.Cells.Item(6, 2).Value = 13.62464540045806
$ws.Cells.Item(6, 3).Value = 6.239178426502917
$ws.Cells.Item(6, 4).Value = 8.01637937001091
$ws.Cells.Item(6, 5).Value = 12.64728385996028
$ws.Cells.Item(6, 6).Value = 36.39415981033878
$ws.Cells.Item(6, 9).Value = 27.29352528300624
$ws.Cells.Item(6, 10).Value = 9.928702792736194
$ws.Cells.Item(6, 11).Value = 10.99744625574507
$ws.Cells.Item(6, 12).Value = 11.16047887015128
$ws.Cells.Item(6, 14).Value = 20.64216712565739
$ws.Cells.Item(6, 15).Value = 28.09789124916143
$ws.Cells.Item(7, 2).Value = 13.68082479454698
$ws.Cells.Item(7, 3).Value = 6.258150972528319
$ws.Cells.Item(7, 4).Value = 8.02030511591782
$ws.Cells.Item(7, 5).Value = 12.64151263955447
$ws.Cells.Item(7, 6).Value = 36.36627426264754
$ws.Cells.Item(7, 9).Value = 27.26540295557485
$ws.Cells.Item(7, 10).Value = 9.923051569258067
$ws.Cells.Item(7, 11).Value = 11.0366798632096
$ws.Cells.Item(7, 12).Value = 11.15955980136425
$ws.Cells.Item(7, 14).Value = 20.62436984232017
$ws.Cells.Item(7, 15).Value = 28.07151958960345
$ws.Cells.Item(8, 2).Value = 13.93102355975657
$ws.Cells.Item(8, 3).Value = 6.339972895055588
$ws.Cells.Item(8, 4).Value = 8.040008133825022
$ws.Cells.Item(8, 5).Value = 12.61932162937047
$ws.Cells.Item(8, 6).Value = 36.25807734114454
$ws.Cells.Item(8, 9).Value = 27.15036781110065
$ws.Cells.Item(8, 10).Value = 9.89960071057156
$ws.Cells.Item(8, 11).Value = 11.21185895822735
$ws.Cells.Item(8, 12).Value = 11.15852971036891
$ws.Cells.Item(8, 14).Value = 20.54952493412351
$ws.Cells.Item(8, 15).Value = 27.96544644021834
$ws.Cells.Item(9, 2).Value = 14.42609578587104
$ws.Cells.Item(9, 3).Value = 6.493131670591159
$ws.Cells.Item(9, 4).Value = 8.08700982268876
$ws.Cells.Item(9, 5).Value = 12.58754101516106
$ws.Cells.Item(9, 6).Value = 36.09937763561175
$ws.Cells.Item(9, 9).Value = 26.9573421680561
$ws.Cells.Item(9, 10).Value = 9.858929069083125
$ws.Cells.Item(9, 11).Value = 11.56030456700009
$ws.Cells.Item(9, 12).Value = 11.16730044604774
$ws.Cells.Item(9, 14).Value = 20.41593708477428
$ws.Cells.Item(9, 15).Value = 27.79449176712521
$ws.Cells.Item(10, 2).Value = 14.78917022560778
$ws.Cells.Item(10, 3).Value = 6.600811010714898
$ws.Cells.Item(10, 4).Value = 8.12645952895482
$ws.Cells.Item(10, 5).Value = 12.57130983079064
$ws.Cells.Item(10, 6).Value = 36.01543850277299
$ws.Cells.Item(10, 9).Value = 26.8353605413174
$ws.Cells.Item(10, 10).Value = 9.832261652052688
$ws.Cells.Item(10, 11).Value = 11.81712152375055
$ws.Cells.Item(10, 12).Value = 11.18024224746049
$ws.Cells.Item(10, 14).Value = 20.3257474441327
$ws.Cells.Item(10, 15).Value = 27.69151239848776
$ws.Cells.Item(11, 2).Value = 14.9535079458369
$ws.Cells.Item(11, 3).Value = 6.648663129478511
$ws.Cells.Item(11, 4).Value = 8.145431851617992
$ws.Cells.Item(11, 5).Value = 12.56546295852039
$ws.Cells.Item(11, 6).Value = 35.9843508993861
$ws.Cells.Item(11, 9).Value = 26.7841720196388
$ws.Cells.Item(11, 10).Value = 9.820822272756855
$ws.Cells.Item(11, 11).Value = 11.93367432795187
$ws.Cells.Item(11, 12).Value = 11.18752162232018
$ws.Cells.Item(11, 14).Value = 20.28642881265336
$ws.Cells.Item(11, 15).Value = 27.64958219251898
$ws.Cells.Item(12, 2).Value = 15.01556370951041
$ws.Cells.Item(12, 3).Value = 6.666614043326947
$ws.Cells.Item(12, 4).Value = 8.152759928950069
$ws.Cells.Item(12, 5).Value = 12.56346909467966
$ws.Cells.Item(12, 6).Value = 35.97359939923606
$ws.Cells.Item(12, 9).Value = 26.76540668667873
$ws.Cells.Item(12, 10).Value = 9.816589528473003
$ws.Cells.Item(12, 11).Value = 11.97773307597781
$ws.Cells.Item(12, 12).Value = 11.19047649041996
$ws.Cells.Item(12, 14).Value = 20.27178439067917
$ws.Cells.Item(12, 15).Value = 27.63441159304887
$ws.Cells.Item(13, 2).Value = 15.00220756894028
$ws.Cells.Item(13, 3).Value = 6.662755645578954
$ws.Cells.Item(13, 4).Value = 8.151175370864024
$ws.Cells.Item(13, 5).Value = 12.56388872728044
$ws.Cells.Item(13, 6).Value = 35.97586953137492
$ws.Cells.Item(13, 9).Value = 26.76942062077684
$ws.Cells.Item(13, 10).Value = 9.817496724091654
$ws.Cells.Item(13, 11).Value = 11.96824826441046
$ws.Cells.Item(13, 12).Value = 11.18983131670954
$ws.Cells.Item(13, 14).Value = 20.27492746329247
$ws.Cells.Item(13, 15).Value = 27.63764737983432
$ws.Cells.Item(14, 2).Value = 14.95861714220897
$ws.Cells.Item(14, 3).Value = 6.650143400711769
$ws.Cells.Item(14, 4).Value = 8.146031880656457
$ws.Cells.Item(14, 5).Value = 12.56529451315829
$ws.Cells.Item(14, 6).Value = 35.98344591060301
$ws.Cells.Item(14, 9).Value = 26.78261578520373
$ws.Cells.Item(14, 10).Value = 9.820472058347764
$ws.Cells.Item(14, 11).Value = 11.93730082934555
$ws.Cells.Item(14, 12).Value = 11.18776075732541
$ws.Cells.Item(14, 14).Value = 20.2852191093982
$ws.Cells.Item(14, 15).Value = 27.64831991644414
$ws.Cells.Item(15, 2).Value = 14.93189226543245
$ws.Cells.Item(15, 3).Value = 6.642395742408075
$ws.Cells.Item(15, 4).Value = 8.142899935161187
$ws.Cells.Item(15, 5).Value = 12.56618425320826
$ws.Cells.Item(15, 6).Value = 35.98821958759667
$ws.Cells.Item(15, 9).Value = 26.79077877364109
$ws.Cells.Item(15, 10).Value = 9.822307431244157
$ws.Cells.Item(15, 11).Value = 11.91833346353739
$ws.Cells.Item(15, 12).Value = 11.18651825110405
$ws.Cells.Item(15, 14).Value = 20.29155487402124
$ws.Cells.Item(15, 15).Value = 27.6549493019847
$ws.Cells.Item(16, 2).Value = 14.77840886285252
$ws.Cells.Item(16, 3).Value = 6.597660394053283
$ws.Cells.Item(16, 4).Value = 8.125239965374215
$ws.Cells.Item(16, 5).Value = 12.57172279697076
$ws.Cells.Item(16, 6).Value = 36.0176129792414
$ws.Cells.Item(16, 9).Value = 26.83879239742099
$ws.Cells.Item(16, 10).Value = 9.833023130348135
$ws.Cells.Item(16, 11).Value = 11.80949574802274
$ws.Cells.Item(16, 12).Value = 11.17979437951198
$ws.Cells.Item(16, 14).Value = 20.32835130830485
$ws.Cells.Item(16, 15).Value = 27.69435160373761
$ws.Cells.Item(17, 2).Value = 14.68399926080079
$ws.Cells.Item(17, 3).Value = 6.569921865461087
$ws.Cells.Item(17, 4).Value = 8.114666297803609
$ws.Cells.Item(17, 5).Value = 12.57551358774886
$ws.Cells.Item(17, 6).Value = 36.03746275569885
$ws.Cells.Item(17, 9).Value = 26.8693489923385
$ws.Cells.Item(17, 10).Value = 9.839773774186011
$ws.Cells.Item(17, 11).Value = 11.74262984626279
$ws.Cells.Item(17, 12).Value = 11.1760247636686
$ws.Cells.Item(17, 14).Value = 20.35136169952443
$ws.Cells.Item(17, 15).Value = 27.71978309667501
$ws.Cells.Item(18, 2).Value = 14.62962281444263
$ws.Cells.Item(18, 3).Value = 6.553861272011765
$ws.Cells.Item(18, 4).Value = 8.108681464106811
$ws.Cells.Item(18, 5).Value = 12.57783864752935
$ws.Cells.Item(18, 6).Value = 36.04954780730773
$ws.Cells.Item(18, 9).Value = 26.88732921228944
$ws.Cells.Item(18, 10).Value = 9.843721700070551
$ws.Cells.Item(18, 11).Value = 11.70414671730715
$ws.Cells.Item(18, 12).Value = 11.17398772759187
$ws.Cells.Item(18, 14).Value = 20.36475759896237
$ws.Cells.Item(18, 15).Value = 27.73487326909864
$ws.Cells.Item(19, 2).Value = 14.61120083330232
$ws.Cells.Item(19, 3).Value = 6.548405417370594
$ws.Cells.Item(19, 4).Value = 8.106671855306342
$ws.Cells.Item(19, 5).Value = 12.57865074732861
$ws.Cells.Item(19, 6).Value = 36.05375431409659
$ws.Cells.Item(19, 9).Value = 26.8934865406277
$ws.Cells.Item(19, 10).Value = 9.845069598247321
$ws.Cells.Item(19, 11).Value = 11.69111409903752
$ws.Cells.Item(19, 12).Value = 11.17332060155973
$ws.Cells.Item(19, 14).Value = 20.36932089331901
$ws.Cells.Item(19, 15).Value = 27.74006198349525
$ws.Cells.Item(20, 2).Value = 14.69405745391417
$ws.Cells.Item(20, 3).Value = 6.572885713388599
$ws.Cells.Item(20, 4).Value = 8.115781886088643
$ws.Cells.Item(20, 5).Value = 12.57509508140081
$ws.Cells.Item(20, 6).Value = 36.03528057671355
$ws.Cells.Item(20, 9).Value = 26.86605428694789
$ws.Cells.Item(20, 10).Value = 9.839048418214176
$ws.Cells.Item(20, 11).Value = 11.74975055807834
$ws.Cells.Item(20, 12).Value = 11.17641248555044
$ws.Cells.Item(20, 14).Value = 20.34889555615381
$ws.Cells.Item(20, 15).Value = 27.71702798370515
$ws.Cells.Item(21, 2).Value = 14.97142590725889
$ws.Cells.Item(21, 3).Value = 6.653852580864802
$ws.Cells.Item(21, 4).Value = 8.147538781344872
$ws.Cells.Item(21, 5).Value = 12.56487562939763
$ws.Cells.Item(21, 6).Value = 35.98119284280457
$ws.Cells.Item(21, 9).Value = 26.77872325483299
$ws.Cells.Item(21, 10).Value = 9.819595443961436
$ws.Cells.Item(21, 11).Value = 11.94639323473145
$ws.Cells.Item(21, 12).Value = 11.1883635635722
$ws.Cells.Item(21, 14).Value = 20.2821895711436
$ws.Cells.Item(21, 15).Value = 27.64516592905897
$ws.Cells.Item(22, 2).Value = 15.15165525765901
$ws.Cells.Item(22, 3).Value = 6.70577723186904
$ws.Cells.Item(22, 4).Value = 8.169129237766279
$ws.Cells.Item(22, 5).Value = 12.55947984903771
$ws.Cells.Item(22, 6).Value = 35.95179267687937
$ws.Cells.Item(22, 9).Value = 26.72525338307703
$ws.Cells.Item(22, 10).Value = 9.807459262639705
$ws.Cells.Item(22, 11).Value = 12.07444338213864
$ws.Cells.Item(22, 12).Value = 11.19732934977714
$ws.Cells.Item(22, 14).Value = 20.24001910613137
$ws.Cells.Item(22, 15).Value = 27.60232376598799
$ws.Cells.Item(23, 2).Value = 15.0555774181768
$ws.Cells.Item(23, 3).Value = 6.678157052358767
$ws.Cells.Item(23, 4).Value = 8.157530896347572
$ws.Cells.Item(23, 5).Value = 12.56224252841806
$ws.Cells.Item(23, 6).Value = 35.96693975284306
$ws.Cells.Item(23, 9).Value = 26.7534612785815
$ws.Cells.Item(23, 10).Value = 9.813883855887173
$ws.Cells.Item(23, 11).Value = 12.00615559000906
$ws.Cells.Item(23, 12).Value = 11.1924390920511
$ws.Cells.Item(23, 14).Value = 20.26239617267886
$ws.Cells.Item(23, 15).Value = 27.62481192294641
$ws.Cells.Item(24, 2).Value = 14.6895104491811
$ws.Cells.Item(24, 3).Value = 6.571546110272412
$ws.Cells.Item(24, 4).Value = 8.115277235201658
$ws.Cells.Item(24, 5).Value = 12.57528383428258
$ws.Cells.Item(24, 6).Value = 36.03626504304081
$ws.Cells.Item(24, 9).Value = 26.86754253733255
$ws.Cells.Item(24, 10).Value = 9.839376143263099
$ws.Cells.Item(24, 11).Value = 11.74653140917477
$ws.Cells.Item(24, 12).Value = 11.1762367908992
$ws.Cells.Item(24, 14).Value = 20.3500099795412
$ws.Cells.Item(24, 15).Value = 27.71827210858228
$ws.Cells.Item(25, 2).Value = 14.29204114661066
$ws.Cells.Item(25, 3).Value = 6.452520881846002
$ws.Cells.Item(25, 4).Value = 8.073417207532547
$ws.Cells.Item(25, 5).Value = 12.59488571269647
$ws.Cells.Item(25, 6).Value = 36.13657824024566
$ws.Cells.Item(25, 9).Value = 27.00607662572158
$ws.Cells.Item(25, 10).Value = 9.869365514692822
$ws.Cells.Item(25, 11).Value = 11.46574080582468
$ws.Cells.Item(25, 12).Value = 11.1637803882916
$ws.Cells.Item(25, 14).Value = 20.45067305953399
$ws.Cells.Item(25, 15).Value = 27.83676911625564
